$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Um ponto = 2 páginas pequeno-médias web (JSP + Angular) sem..."
#    -> "Um ponto = 1 páginas pequeno-médias web (JSP + Angular) sem..."
#    Toggle Bold on/off around the Text assignment so the engine keeps
#    this single-character run distinct from its neighbours instead of
#    silently coalescing it into the adjoining identically-formatted
#    run.
# ------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Um ponto = 2 páginas pequeno-médias web", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $rng1.Start
$digitRng = $d.Range($start1 + 11, $start1 + 12)
$digitRng.Font.Bold = $true
$digitRng.Text = "1"
$digitRng2 = $d.Range($start1 + 11, $start1 + 12)
$digitRng2.Font.Bold = $false

# ------------------------------------------------------------------
# 2) "Um ponto = 1 páginas pequeno-média web (JSP + Angular) e 1 Servlets
#     de comunicação com o back-end."
#    -> "Um ponto = 1 Servlets de comunicação com o back-end."
#    A pure Range.Delete() leaves untouched run boundaries alone, so we
#    use it to drop the middle portion of the sentence (this also
#    removes the old "_GoBack" bookmark that used to sit right before
#    the trailing space / "Servlets").
# ------------------------------------------------------------------
$rUm = $d.Content
$rUm.Find.Execute("Um ponto = 1 páginas pequeno-média web", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$umStart = $rUm.Start
$delStart = $umStart + 12

$rServlets = $d.Content
$rServlets.Find.Execute("Servlets de comunicação com o back-end.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$servletsStart = $rServlets.Start
$delEnd = $servletsStart - 1

$delRng = $d.Range($delStart, $delEnd)
$delRng.Delete()

# ------------------------------------------------------------------
# 3) Re-seat the "_GoBack" bookmark at the last edited spot: right
#    after the "3" in "...trabalha em média 3 horas por dia...".
#    Bookmarks.Add("_GoBack", ...) removes any pre-existing _GoBack
#    bookmark automatically, which is exactly what we need here.
# ------------------------------------------------------------------
$rHoras = $d.Content
$rHoras.Find.Execute("trabalha em média 3 horas", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$horasStart = $rHoras.Start
$bmPos = $horasStart + 19
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)
